$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Update row 5: Bellamy -> Yijia, email change
$ws.Range("B5").Value = "Yijia"
$ws.Range("D5").Value = "bellamy93158@gmail.com"

# Add new rows 6-8
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Yijia"
$ws.Range("C6").Value = "Sun"
$ws.Range("D6").Value = "michaelice2604@gmail.com"
$ws.Range("E6").Value = "Shandong"
$ws.Range("F6").Value = 8

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Yijia"
$ws.Range("C7").Value = "Sun"
$ws.Range("D7").Value = "michaelice2604@gmail.com"
$ws.Range("E7").Value = "Florida"
$ws.Range("F7").Value = -5

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Yijia"
$ws.Range("C8").Value = "Sun"
$ws.Range("D8").Value = "michaelice2604@gmail.com"
$ws.Range("E8").Value = "Hongkong"
$ws.Range("F8").Value = "8"

# Apply bold/centered style to column A for new rows, matching existing header style
$ws.Range("A6:A8").Font.Bold = $true
$ws.Range("A6:A8").Font.Name = "等线"
$ws.Range("A6:A8").Font.Size = 11
$ws.Range("A6:A8").HorizontalAlignment = -4108
$ws.Range("A6:A8").VerticalAlignment = -4160
$ws.Range("A6:A8").Borders.LineStyle = 1

# selection
$ws.Range("G12").Select()
